# "fix log in/out mechanism"
# - Renames the "View Question Table..." feature row to "View Response Table..."
#   (login/logout now lives on the contact-form "response" table, not a
#   question table).
# - Normalizes a couple of stray/duplicate cell styles on the Admin sheet
#   (no visual change — B8 was using a fill/border xf that is identical to
#   the plain date style used by its neighbours).
# - Extends the category-header merges on Admin to span column C, matching
#   the other header rows.
# - Adds a new "Scouting Admin" sheet (right after "Admin") that tracks the
#   same kind of feature/worked/notes rows for the Scouting side of the
#   login fix, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$admin = $wb.Worksheets.Item("Admin")

# ---------------------------------------------------------------------
# 1) Admin sheet: rename the feature description text.
# ---------------------------------------------------------------------
$admin.Range("A14").Value = "View Response Table with row click action"
$admin.Range("A18").Value = "View Response Table with row click action"

# ---------------------------------------------------------------------
# 2) Admin sheet: tidy up B8's style so it matches the other plain date
#    cells (B3:B6) instead of the redundant fill/border variant.
# ---------------------------------------------------------------------
$admin.Range("B3").Copy()
$admin.Range("B8").PasteSpecial(-4122)
$admin.Range("B8").Value = 45314

# ---------------------------------------------------------------------
# 3) Admin sheet: widen the section-header merges to include column C.
#    (A2:C2 is already merged in the source file, so only the B-only
#    merges need to be extended.)
# ---------------------------------------------------------------------
$admin.Range("A10:C10").Merge()
$admin.Range("A11:C11").Merge()
$admin.Range("A15:C15").Merge()

$admin.Range("A1:C6").Select()

# ---------------------------------------------------------------------
# 4) Add the new "Scouting Admin" sheet right after "Admin".
# ---------------------------------------------------------------------
$scouting = $wb.Worksheets.Add($null, $admin)
$scouting.Name = "Scouting Admin"

$scouting.Range("A1").Value = "Feature"
$scouting.Range("B1").Value = "Works"
$scouting.Range("C1").Value = "Notes"

$scouting.Range("A2").Value = "Manage Users"

$scouting.Range("A3").Value = "User Table Filter"
$scouting.Range("B3").Value = 45315

$scouting.Range("A4").Value = "User Table Main Display Edit Capabilities"
$scouting.Range("B4").Value = 45314

$scouting.Range("A5").Value = "Manage Users Modal"
$scouting.Range("B5").Value = 45314

# Match formatting from the Admin sheet: header/category-row style on row 2,
# and the plain date style on the B column date cells.
$admin.Range("A2:C2").Copy()
$scouting.Range("A2:C2").PasteSpecial(-4122)
$scouting.Range("A2:C2").Merge()

$admin.Range("B3").Copy()
$scouting.Range("B3:B5").PasteSpecial(-4122)
$scouting.Range("B3").Value = 45315
$scouting.Range("B4").Value = 45314
$scouting.Range("B5").Value = 45314

$scouting.Range("A4:XFD4").Select()
$scouting.Activate()
